# Scheduled-runner refresh of market-derived Leve profit columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, columns H-N)
# across all eight crafter sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1289.25
$ws.Range("J17").Value = 1321.8667
$ws.Range("L17").Value = 3965.6001
$ws.Range("N17").Value = -4301.6001

$ws.Range("H43").Value = 211100.9
$ws.Range("I43").Value = 7800
$ws.Range("J43").Value = 414401.8
$ws.Range("K43").Value = 7800
$ws.Range("L43").Value = 414401.8
$ws.Range("M43").Value = -7731
$ws.Range("N43").Value = -414539.8

$ws.Range("H64").Value = 4727.091
$ws.Range("I64").Value = 3499
$ws.Range("K64").Value = 3499
$ws.Range("M64").Value = -3251

$ws.Range("H67").Value = 4727.091
$ws.Range("I67").Value = 3499
$ws.Range("K67").Value = 3499
$ws.Range("M67").Value = -2641

$ws.Range("H137").Value = 3722.25
$ws.Range("I137").Value = 2617.158
$ws.Range("J137").Value = 6055.222
$ws.Range("K137").Value = 7851.474
$ws.Range("L137").Value = 18165.666
$ws.Range("M137").Value = -5301.474
$ws.Range("N137").Value = -23265.666

$ws.Range("H138").Value = 2488.4343
$ws.Range("I138").Value = 1365.36
$ws.Range("J138").Value = 3038.9607
$ws.Range("K138").Value = 4096.08
$ws.Range("L138").Value = 9116.882100000001
$ws.Range("M138").Value = 1043.92
$ws.Range("N138").Value = -19396.8821

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8200066
$ws.Range("I32").Value = 8931350
$ws.Range("J32").Value = 9686
$ws.Range("K32").Value = 8931350
$ws.Range("L32").Value = 9686
$ws.Range("M32").Value = -8931063
$ws.Range("N32").Value = -10260

$ws.Range("H61").Value = 11140414
$ws.Range("I61").Value = 17245282
$ws.Range("K61").Value = 17245282
$ws.Range("M61").Value = -17245070

$ws.Range("H64").Value = 49332.668
$ws.Range("J64").Value = 49332.668
$ws.Range("L64").Value = 49332.668
$ws.Range("N64").Value = -49828.668

$ws.Range("H67").Value = 49332.668
$ws.Range("J67").Value = 49332.668
$ws.Range("L67").Value = 49332.668
$ws.Range("N67").Value = -51048.668

$ws.Range("H112").Value = 94995
$ws.Range("J112").Value = 94995
$ws.Range("L112").Value = 94995
$ws.Range("N112").Value = -97949

$ws.Range("H132").Value = 3404.0605
$ws.Range("I132").Value = 1555.1578
$ws.Range("J132").Value = 5913.2856
$ws.Range("K132").Value = 4665.4734
$ws.Range("L132").Value = 17739.8568
$ws.Range("M132").Value = -2135.4734
$ws.Range("N132").Value = -22799.8568

$ws.Range("H136").Value = 11140414
$ws.Range("I136").Value = 17245282
$ws.Range("K136").Value = 51735846
$ws.Range("M136").Value = -51733296

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 886.55554
$ws.Range("I22").Value = 998.4286
$ws.Range("K22").Value = 998.4286
$ws.Range("M22").Value = -825.4286

$ws.Range("H94").Value = 939.5714
$ws.Range("J94").Value = 877.1111
$ws.Range("L94").Value = 877.1111
$ws.Range("N94").Value = -1779.1111

$ws.Range("H126").Value = 34999
$ws.Range("J126").Value = 34999
$ws.Range("L126").Value = 34999
$ws.Range("N126").Value = -44879

$ws.Range("H134").Value = 419034.53
$ws.Range("I134").Value = 2114.8096
$ws.Range("K134").Value = 6344.4288
$ws.Range("M134").Value = -3809.4288

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 537890.9399999999
$ws.Range("I31").Value = 10708.55
$ws.Range("J31").Value = 977209.5600000001
$ws.Range("K31").Value = 10708.55
$ws.Range("L31").Value = 977209.5600000001
$ws.Range("M31").Value = -10413.55
$ws.Range("N31").Value = -977799.5600000001

$ws.Range("H34").Value = 537890.9399999999
$ws.Range("I34").Value = 10708.55
$ws.Range("J34").Value = 977209.5600000001
$ws.Range("K34").Value = 10708.55
$ws.Range("L34").Value = 977209.5600000001
$ws.Range("M34").Value = -10506.55
$ws.Range("N34").Value = -977613.5600000001

$ws.Range("H58").Value = 1420.5
$ws.Range("I58").Value = 1420.5
$ws.Range("K58").Value = 1420.5
$ws.Range("M58").Value = -1217.5

$ws.Range("H62").Value = 2627.4285
$ws.Range("I62").Value = 2518.6
$ws.Range("J62").Value = 2899.5
$ws.Range("K62").Value = 2518.6
$ws.Range("L62").Value = 2899.5
$ws.Range("M62").Value = -1894.6
$ws.Range("N62").Value = -4147.5

$ws.Range("H65").Value = 2627.4285
$ws.Range("I65").Value = 2518.6
$ws.Range("J65").Value = 2899.5
$ws.Range("K65").Value = 12593
$ws.Range("L65").Value = 14497.5
$ws.Range("M65").Value = -9473
$ws.Range("N65").Value = -20737.5

$ws.Range("H86").Value = 6400
$ws.Range("I86").Value = 6850
$ws.Range("K86").Value = 6850
$ws.Range("M86").Value = -5727

$ws.Range("H89").Value = 6400
$ws.Range("I89").Value = 6850
$ws.Range("K89").Value = 34250
$ws.Range("M89").Value = -28634

$ws.Range("H132").Value = 2348
$ws.Range("I132").Value = 2144.9546
$ws.Range("J132").Value = 3837
$ws.Range("K132").Value = 6434.8638
$ws.Range("L132").Value = 11511
$ws.Range("M132").Value = -3904.8638
$ws.Range("N132").Value = -16571

$ws.Range("H134").Value = 2502024.5
$ws.Range("I134").Value = 2502024.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7506073.5
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -7503538.5

$ws.Range("H136").Value = 1420.5
$ws.Range("I136").Value = 1420.5
$ws.Range("K136").Value = 4261.5
$ws.Range("M136").Value = -1711.5

$ws.Range("H141").Value = 216299
$ws.Range("J141").Value = 216299
$ws.Range("L141").Value = 216299
$ws.Range("N141").Value = -226659

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 34333332
$ws.Range("I4").Value = 101000000
$ws.Range("J4").Value = 1000000
$ws.Range("K4").Value = 303000000
$ws.Range("L4").Value = 3000000
$ws.Range("M4").Value = -302999888
$ws.Range("N4").Value = -3000224

$ws.Range("H5").Value = 1701.9445
$ws.Range("J5").Value = 1757.375
$ws.Range("L5").Value = 5272.125
$ws.Range("N5").Value = -5496.125

$ws.Range("H92").Value = 627008.5
$ws.Range("I92").Value = 1112952
$ws.Range("J92").Value = 2224
$ws.Range("K92").Value = 3338856
$ws.Range("L92").Value = 6672
$ws.Range("M92").Value = -3337608
$ws.Range("N92").Value = -9168

$ws.Range("H124").Value = 1682.7273
$ws.Range("I124").Value = 1682.7273
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 5048.1819
$ws.Range("L124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -138.1818999999996

$ws.Range("H125").Value = 8277.214
$ws.Range("J125").Value = 10003
$ws.Range("L125").Value = 30009
$ws.Range("N125").Value = -39849

$ws.Range("H129").Value = 47620268
$ws.Range("J129").Value = 83335260
$ws.Range("L129").Value = 250005780
$ws.Range("N129").Value = -250015780

$ws.Range("H135").Value = 1701.9445
$ws.Range("J135").Value = 1757.375
$ws.Range("L135").Value = 15816.375
$ws.Range("N135").Value = -20886.375

$ws.Range("H140").Value = 232343.31
$ws.Range("I140").Value = 232343.31
$ws.Range("K140").Value = 697029.9299999999
$ws.Range("M140").Value = -691849.9299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3165.9167
$ws.Range("I80").Value = 2888
$ws.Range("K80").Value = 2888
$ws.Range("M80").Value = -1890

$ws.Range("H83").Value = 3165.9167
$ws.Range("I83").Value = 2888
$ws.Range("K83").Value = 14440
$ws.Range("M83").Value = -9448

$ws.Range("H111").Value = 130011.336
$ws.Range("J111").Value = 130011.336
$ws.Range("L111").Value = 130011.336
$ws.Range("N111").Value = -136145.336

$ws.Range("H132").Value = 62508960
$ws.Range("I132").Value = 83334990
$ws.Range("J132").Value = 30853.25
$ws.Range("K132").Value = 250004970
$ws.Range("L132").Value = 92559.75
$ws.Range("M132").Value = -250002440
$ws.Range("N132").Value = -97619.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 974.2308
$ws.Range("I16").Value = 830.7727
$ws.Range("J16").Value = 1763.25
$ws.Range("K16").Value = 830.7727
$ws.Range("L16").Value = 1763.25
$ws.Range("M16").Value = -660.7727
$ws.Range("N16").Value = -2103.25

$ws.Range("H136").Value = 48484.703
$ws.Range("I136").Value = 6673.722
$ws.Range("J136").Value = 132106.67
$ws.Range("K136").Value = 20021.166
$ws.Range("L136").Value = 396320.01
$ws.Range("M136").Value = -17471.166
$ws.Range("N136").Value = -401420.01

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 22263.5
$ws.Range("J41").Value = 22263.5
$ws.Range("L41").Value = 22263.5
$ws.Range("N41").Value = -23043.5

$ws.Range("H75").Value = 99980
$ws.Range("J75").Value = 99980
$ws.Range("L75").Value = 99980
$ws.Range("N75").Value = -101852

$ws.Range("H78").Value = 99980
$ws.Range("J78").Value = 99980
$ws.Range("L78").Value = 299940
$ws.Range("N78").Value = -309300

$ws.Range("H81").Value = 3150
$ws.Range("I81").Value = 3150
$ws.Range("K81").Value = 6300
$ws.Range("M81").Value = -5239

$ws.Range("H84").Value = 3150
$ws.Range("I84").Value = 3150
$ws.Range("K84").Value = 31500
$ws.Range("M84").Value = -26196

$ws.Range("H98").Value = 24012.715
$ws.Range("J98").Value = 24012.715
$ws.Range("L98").Value = 24012.715
$ws.Range("N98").Value = -30002.715
